# Update "Latest HO Xliff Generate/Handoff/Handback" timestamps for the
# c7e278f5.../6c2dd570... entries to reflect a freshly regenerated report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
# Row 3 = 6c2dd570-ba5c-49fb-8606-9313d1334d4f, column G = "Latest HO Xliff
# Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-24 04:43:31"

# --- zh-cn sheet -------------------------------------------------------------
# Row 3 = 6c2dd570-ba5c-49fb-8606-9313d1334d4f
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-24 04:43:26"
$wsZhCn.Range("K3").Value = "2016-08-24 04:43:53"

# --- de-de sheet -------------------------------------------------------------
# Row 3 = 6c2dd570-ba5c-49fb-8606-9313d1334d4f
# Column H = "Correspond Handoff Datetime" (shares text with Overview!G3)
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-24 04:43:31"
$wsDeDe.Range("K3").Value = "2016-08-24 04:43:59"
